$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 296 (shifts old rows 296-315 down to 297-316,
# extending the used range to A1:R316).
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new weekly price record.
$ws.Cells.Item(296, 1).Value = 4
$ws.Cells.Item(296, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(296, 3).Value = "Los Lagos"
$ws.Cells.Item(296, 4).Value = 44783
$ws.Cells.Item(296, 5).Value = 10
$ws.Cells.Item(296, 6).Value = 100112037
$ws.Cells.Item(296, 7).Value = "Cebollín"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Segunda"
$ws.Cells.Item(296, 10).Value = 35
$ws.Cells.Item(296, 11).Value = 9500
$ws.Cells.Item(296, 12).Value = 9500
$ws.Cells.Item(296, 13).Value = 9500
$ws.Cells.Item(296, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(296, 15).Value = "Región Metropolitana"
$ws.Cells.Item(296, 16).Value = 264
$ws.Cells.Item(296, 17).Value = 36
$ws.Cells.Item(296, 18).Value = "Hortaliza"
